# Fruta / hortaliza, semanal
#
# The weekly refresh reshuffles which source record (date, volume,
# min/max/weighted price, unit, origin, price-per-kg, kg-per-unit) lands on
# each existing row of the sheet. The row-identity columns (Mercado ID,
# Mercado, Region, Codreg, Tipo, Producto ID, Producto, Categoria ID,
# Categoria, Variedad, Calidad) are identical for every data row and stay
# untouched; only columns D and M:T move between rows according to the
# mapping below (destination row -> source row, both referring to the
# *original* sheet state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2  = 19
    3  = 21
    4  = 14
    5  = 5
    6  = 9
    7  = 13
    8  = 16
    9  = 11
    10 = 2
    11 = 20
    12 = 10
    13 = 6
    14 = 4
    15 = 18
    16 = 7
    17 = 15
    18 = 3
    19 = 17
    20 = 12
    21 = 8
}

$cols = @("D", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the original values for the columns that move, before any writes,
# so chained writes don't clobber a row that still needs to be read as a
# source for another destination row.
$snapshot = @{}
foreach ($row in $rowMap.Keys) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value()
    }
    $snapshot[$row] = $rowValues
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcValues[$col]
    }
}
